$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.635.96'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.900.68'
$ws.Range("E3").Value = '  -4.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.81'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.31'
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("E8").Value = '  -3.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.900.47'
$ws.Range("E9").Value = '  -3.99%  '
$ws.Range("E10").Value = '  +5.39%  '
$ws.Range("E11").Value = '  -4.63%  '
$ws.Range("E12").Value = '  -2.62%  '
$ws.Range("E13").Value = '  -4.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.76'
$ws.Range("E14").Value = '  -2.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.378.43'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.573.64'
$ws.Range("E17").Value = '  -3.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.79'
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.898.27'
$ws.Range("E19").Value = '  -4.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '423.97'
$ws.Range("E20").Value = '  -5.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("E21").Value = '  -4.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.669'
$ws.Range("E22").Value = '  -3.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.08'
$ws.Range("E23").Value = '  -5.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.21'
$ws.Range("E24").Value = '  -2.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.02'
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.86'
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.20'
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("E32").Value = '  -3.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.38'
$ws.Range("E33").Value = '  -4.18%  '
$ws.Range("E34").Value = '  -2.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  -2.58%  '
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.63'
$ws.Range("E37").Value = '  -3.72%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.96'
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.14'
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("E40").Value = '  -2.60%  '
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.69'
$ws.Range("E42").Value = '  -3.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.291'
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.47'
$ws.Range("E44").Value = '  +2.89%  '
$ws.Range("E45").Value = '  -2.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '370.18'
$ws.Range("E46").Value = '  -5.37%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.657.36'
$ws.Range("E47").Value = '  -2.90%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.29'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.39'
$ws.Range("E49").Value = '  +6.53%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("E51").Value = '  -1.41%  '
